# "Depositar" functionality: record a set of deposit/expense entries on the
# "Sheet" tab and extend the running balance on the "Saldo" tab to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet")
$ws2 = $wb.Worksheets.Item("Saldo")

# --- Sheet: replace the two sample rows with the real deposit history ---
$ws1.Range("A2:E3").ClearContents()

# Row 2: initial deposit of 100
$ws1.Range("A2").Value = 100
$ws1.Range("C2").Value = "26/02/2023"

# Rows 3-5: further deposits recorded as plain text amounts
$ws1.Range("A3:A6").NumberFormat = "@"

$ws1.Range("A3").Value = "56.36"
$ws1.Range("C3").Value = "26/02/2023"

$ws1.Range("A4").Value = "85.90"
$ws1.Range("C4").Value = "26/02/2023"

$ws1.Range("A5").Value = "1.58"
$ws1.Range("C5").Value = "26/02/2023"

# Row 6: a withdrawal/expense entry (haircut) referencing account 123/321
$ws1.Range("A6").Value = "25.63"
$ws1.Range("B6").Value = "123/321"
$ws1.Range("C6").Value = "26/02/2023"
$ws1.Range("D6").Value = "cabelo"

# --- Saldo: extend the running balance to reflect the new entries ---
$ws2.Range("A3").Value = 100
$ws2.Range("A4").Value = 156.36
$ws2.Range("A5").Value = 242.26
$ws2.Range("A6").Value = 243.84
$ws2.Range("A7").Value = 218.21

# --- Selection / active sheet bookkeeping ---
[void]$ws2.Activate()
[void]$ws2.Range("A1").Select()
[void]$ws1.Activate()
[void]$ws1.Range("A1").Select()
